$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus indices 0-23), columns B-F and I-N
# (column G stays 1, column H has no data, column A is the bus index - unchanged)
$newData = @(
    @(1.02, 1.019008834503968, 1.024525588268876, 1.020190472673558, 1.032940651327546, 1.029328067124957, 1.024214032409659, 1.027353833164618, 1.023031491326732, 1.035744433281276, 1.012108334454404),
    @(1.02, 1.020123686371282, 1.025350889887488, 1.021140823844667, 1.035186782754046, 1.029632716461576, 1.024964587939291, 1.027986195009908, 1.023787632899246, 1.03779560589015, 1.012361262998139),
    @(1.02, 1.02084458527579, 1.025884109190973, 1.021755709313386, 1.036633551855234, 1.029827392179558, 1.025449262252439, 1.028393928027014, 1.02427623824209, 1.039115915420564, 1.012524466671664),
    @(1.02, 1.021147537591646, 1.026108083243806, 1.022014195338225, 1.037240222885338, 1.029908649152297, 1.025652785139086, 1.02856499422881, 1.024481489259636, 1.039669346633664, 1.012592968350483),
    @(1.02, 1.021198398029371, 1.026145678298791, 1.022057595626102, 1.037341995725342, 1.029922258340247, 1.025686943842487, 1.02859369683006, 1.024515942545751, 1.039762175734758, 1.01260446369783),
    @(1.02, 1.020848633782028, 1.025887102692573, 1.021759163260527, 1.03664166426801, 1.029828480233165, 1.025451982652866, 1.028396215174026, 1.02427898143745, 1.039123316743017, 1.012525382422532),
    @(1.02, 1.019385705795652, 1.02480467069812, 1.020511660606665, 1.033701136233717, 1.029431534177024, 1.024467891204816, 1.027567843564385, 1.023287172274269, 1.036439093278321, 1.012193907978995),
    @(1.02, 1.016804013947487, 1.022891038978034, 1.018312913156336, 1.028467149423053, 1.028713164883348, 1.022726172640136, 1.026096981365483, 1.021534288975039, 1.031654497740059, 1.011606276739863),
    @(1.02, 1.01508015174329, 1.021610973705911, 1.016846658024074, 1.024940281423109, 1.02822138133988, 1.021559788894972, 1.025108782824719, 1.020362114640171, 1.028425866110168, 1.011212114801513),
    @(1.02, 1.014333014472, 1.021055642916056, 1.016211630432064, 1.02340368262552, 1.028005343314555, 1.021053463261541, 1.024679046449133, 1.019853677087561, 1.027018119693067, 1.011040859274587),
    @(1.02, 1.014055386873034, 1.02084920759535, 1.015975731400945, 1.02283145759049, 1.027924629206382, 1.020865197627453, 1.024519144326626, 1.019664686712495, 1.026493716215581, 1.010977159354904),
    @(1.02, 1.014114943856643, 1.020893495984955, 1.016026333550788, 1.02295426878061, 1.027941963887619, 1.020905590051334, 1.0245534565358, 1.019705231888318, 1.026606271278455, 1.01099082719765),
    @(1.02, 1.014310067901615, 1.021038582201037, 1.016192131395422, 1.023356412452446, 1.027998681027222, 1.021037905131868, 1.024665834593779, 1.019838057830502, 1.026974803262328, 1.011035595618654),
    @(1.02, 1.014430275913095, 1.021127953235617, 1.016294282006955, 1.023603991099667, 1.028033564217207, 1.021119403074221, 1.024735037392428, 1.01991987845436, 1.027201667435825, 1.011063167238431),
    @(1.02, 1.015129721885994, 1.021647806771811, 1.016888799811779, 1.025042057650366, 1.028235653632051, 1.021593365020736, 1.025137264038259, 1.020395839284858, 1.028519085060233, 1.011223468167754),
    @(1.02, 1.015568277829754, 1.021973613294515, 1.017261689026247, 1.025941559271135, 1.028361588686789, 1.021890326136903, 1.02538907586546, 1.020694160336274, 1.029342832348577, 1.011323864695111),
    @(1.02, 1.015824013403481, 1.022163549188213, 1.017479176649505, 1.026465315466092, 1.028434746322201, 1.022063415795117, 1.025535776092644, 1.020868081355248, 1.029822374528934, 1.011382368302168),
    @(1.02, 1.015911201413948, 1.022228295267419, 1.01755333227569, 1.026643750159918, 1.028459640713858, 1.022122414148178, 1.025585767070941, 1.020927369639684, 1.029985728581814, 1.011402307045419),
    @(1.02, 1.015521231849691, 1.021938667857876, 1.017221682834605, 1.025845145492011, 1.028348107910583, 1.02185847773204, 1.02536207720172, 1.020662162076676, 1.029254549145687, 1.011313098894368),
    @(1.02, 1.014252611729779, 1.02099586238314, 1.016143308673105, 1.023238031941127, 1.027981992188688, 1.020998947006301, 1.024632749775784, 1.019798947617849, 1.026866321633886, 1.011022414869456),
    @(1.02, 1.013454355960376, 1.020402151953562, 1.015465166539779, 1.021590343868565, 1.027749091480602, 1.020457403734944, 1.024172577787411, 1.019255433914406, 1.025356026269213, 1.010839140817871),
    @(1.02, 1.013877586721547, 1.020716978161476, 1.015824675110585, 1.022464634817186, 1.027872814493376, 1.020744593246387, 1.024416677578936, 1.019543635118603, 1.026157503493056, 1.010936346400376),
    @(1.02, 1.01554249010948, 1.021954458510562, 1.017239759949103, 1.025888713532147, 1.028354200214374, 1.021872869034841, 1.025374277283963, 1.020676620974563, 1.029294443414412, 1.011317963669427),
    @(1.02, 1.017471913078117, 1.023386510215152, 1.018881409605203, 1.029826705752743, 1.028901136892453, 1.023177363495915, 1.026478569510703, 1.021988076057855, 1.032898126495004, 1.011758615361263)
)

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$startRow = 2
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $rowVals[$j]
    }
}
